$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header addition
$ws.Range("S1").Value = "Actual Spend Total"

# Row 2
$ws.Range("A2").Value = "fixedPrice/2023-08-24T08:30:26.178Z/client2"
$ws.Range("B2").Value = "newProject"
$ws.Range("C2").Formula = "="""""
$ws.Range("D2").Value = "fixedPrice"
$ws.Range("E2").Value = "notStarted"
$ws.Range("F2").Value = 45162.35446965278
$ws.Range("G2").Value = "Ahmed shalaab"
$ws.Range("H2").Value = "client2pm"
$ws.Range("K2").Value = "USD"
$ws.Range("L2").Value = 35234
$ws.Range("M2").Value = "Signed"
$ws.Range("N2").Value = "referenceNumber1"
$ws.Range("O2").Value = 12345
$ws.Range("Q2").Value = "client2"
$ws.Range("S2").Value = 4000

# Row 3
$ws.Range("A3").Value = "TnM/3122-06-12T11:31:00Z/client1"
$ws.Range("B3").Value = "iprojectnase"
$ws.Range("C3").Formula = "="""""
$ws.Range("D3").Value = "TnM"
$ws.Range("E3").Value = "notStarted"
$ws.Range("F3").Value = 446490.4798611111
$ws.Range("G3").Value = "obaid saafan"
$ws.Range("H3").Value = "client1pm"
$ws.Range("K3").Value = "USD"
$ws.Range("L3").Value = 35234
$ws.Range("M3").Value = "Signed"
$ws.Range("N3").Value = "referenceNumber1"
$ws.Range("O3").Value = 12345
$ws.Range("P3").Value = 12345
$ws.Range("Q3").Value = "client1"
$ws.Range("S3").Value = 56000

# Both project-start-date cells use the same date number format (numFmtId 14).
# Set it once on F2, then copy just the format over to F3 so they share one style.
$ws.Range("F2").NumberFormat = "mm-dd-yy"
$ws.Range("F2").Copy()
$ws.Range("F3").PasteSpecial(-4122)
$excel.CutCopyMode = $false
